$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for 2024-03-04
$ws.Range("A10").Value = "'20240304"
$ws.Range("B10").Value = 109
$ws.Range("C10").Value = 265
$ws.Range("D10").Value = 84
$ws.Range("E10").Value = 133
$ws.Range("F10").Value = 214
$ws.Range("G10").Value = 77

# New row for 2024-03-05
$ws.Range("A11").Value = "'20240305"
$ws.Range("B11").Value = 193
$ws.Range("C11").Value = 449
$ws.Range("D11").Value = 69
$ws.Range("E11").Value = 123
$ws.Range("F11").Value = 248
$ws.Range("G11").Value = 72

# Match the Date column formatting (bold/border/centered) used by the rows above
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
